$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06686833333333334
$ws.Range("N2").Value = 0.200605
$ws.Range("O2").Value = 0.05567274787007094
$ws.Range("P2").Value = 0.05567274787007093
$ws.Range("Q2").Value = 10.76395318074945
$ws.Range("R2").Value = 96.87557862674501
$ws.Range("S2").Value = 0.02188855172981137
$ws.Range("T2").Value = 0.02188855172981137
$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.8610905203551533
$ws.Range("P3").Value = 0.8610905203551533
$ws.Range("Q3").Value = 166.4860887973669
$ws.Range("R3").Value = 1498.374799176302
$ws.Range("S3").Value = 0.3385502803423945
$ws.Range("T3").Value = 0.3385502803423945
$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.08323673177477579
$ws.Range("P4").Value = 0.08323673177477578
$ws.Range("Q4").Value = 16.09326498187712
$ws.Range("R4").Value = 144.839384836894
$ws.Range("S4").Value = 0.03272573348677952
$ws.Range("T4").Value = 0.03272573348677951
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06686833333333334
$ws.Range("N5").Value = 0.200605
$ws.Range("O5").Value = 0.05567274787007094
$ws.Range("P5").Value = 0.05567274787007093
$ws.Range("Q5").Value = 6.016626627919444
$ws.Range("R5").Value = 54.149639651275
$ws.Range("S5").Value = 0.01223483983744028
$ws.Range("T5").Value = 0.01223483983744028
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.8610905203551533
$ws.Range("P6").Value = 0.8610905203551533
$ws.Range("S6").Value = 0.1892362961259017
$ws.Range("T6").Value = 0.1892362961259017
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.08323673177477579
$ws.Range("P7").Value = 0.08323673177477578
$ws.Range("S7").Value = 0.01829239836038041
$ws.Range("T7").Value = 0.01829239836038041
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.06686833333333334
$ws.Range("N8").Value = 0.200605
$ws.Range("O8").Value = 0.05567274787007094
$ws.Range("P8").Value = 0.05567274787007093
$ws.Range("Q8").Value = 10.59714983348667
$ws.Range("R8").Value = 95.37434850138
$ws.Range("S8").Value = 0.02154935630281929
$ws.Range("T8").Value = 0.0215493563028193
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.8610905203551533
$ws.Range("P9").Value = 0.8610905203551533
$ws.Range("S9").Value = 0.3333039438868571
$ws.Range("T9").Value = 0.3333039438868572
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.08323673177477579
$ws.Range("P10").Value = 0.08323673177477578
